$wb = $excel.ActiveWorkbook

# New row of data (row 43) to append to each of the 4 sheets, matching the
# existing "time / 总长 / ID / 实际长度 / 和校验 / ..._DEC" table layout.
$rows = @{
    "FE_LFT_#1" = @{
        A = 45829.49726851852
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x60"
        E = "0xf"
        F = 380
        G = [double]"7.598631275147109e+23"
        H = 352
        I = 15
    }
    "FE_LFT_#2" = @{
        A = 45829.49726851852
        B = "0x01,0x90"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x74"
        E = "0xe"
        F = 400
        G = [double]"5.68432987514711e+23"
        H = 372
        I = 14
    }
    "FE_PLT_#1" = @{
        A = 45829.49726851852
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x6A"
        E = "0x3"
        F = 110
        G = [double]"5.68631262647114e+23"
        H = 106
        I = 3
    }
    "FE_PLT_#2" = @{
        A = 45829.49726851852
        B = "0x00,0x6e"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x6A"
        E = "0x3"
        F = 110
        G = [double]"9.85046333984776e+23"
        H = 106
        I = 3
    }
}

foreach ($sheetName in $rows.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $data = $rows[$sheetName]

    $ws.Range("A43").Value = $data.A
    $ws.Range("A43").NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Range("B43").Value = $data.B
    $ws.Range("C43").Value = $data.C
    $ws.Range("D43").Value = $data.D
    $ws.Range("E43").Value = $data.E

    $ws.Range("F43").Value = $data.F
    $ws.Range("G43").Value = $data.G
    $ws.Range("H43").Value = $data.H
    $ws.Range("I43").Value = $data.I
}
